# UPDATE truc a alex + feuille de temps
#
# Applies the "10-05" (sheet3) timesheet updates:
#  - row 25: task description changed ("DCO + DS" -> "DCO + DS + DN") and the
#    time entry bumped from 2 to 4 hours.
#  - rows 41-49: new timesheet entries filled in (previously blank rows).
#  - sheet view: selection/scroll moved down to the newly-filled-in area.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10-05")

# --- Row 25: bump the existing entry's hours from 2 to 4 -------------------
$ws.Cells.Item(25, 6).Value = 4

# --- Fill in date / member / category / hours for the new rows first -------
# (these columns only ever reference members & category codes that already
# exist in the shared-string table, so the order they're written in here
# doesn't affect the resulting workbook.)

# Row 41: Gabriel / DV / 3h
$ws.Cells.Item(41, 2).Value = 42284
$ws.Cells.Item(41, 3).Value = "Gabriel"
$ws.Cells.Item(41, 4).Value = "DV"
$ws.Cells.Item(41, 6).Value = 3

# Row 42: David / DV / 3h
$ws.Cells.Item(42, 2).Value = 42284
$ws.Cells.Item(42, 3).Value = "David"
$ws.Cells.Item(42, 4).Value = "DV"
$ws.Cells.Item(42, 6).Value = 3

# Row 43: Alex / DC / 3h
$ws.Cells.Item(43, 2).Value = 42284
$ws.Cells.Item(43, 3).Value = "Alex"
$ws.Cells.Item(43, 4).Value = "DC"
$ws.Cells.Item(43, 6).Value = 3

# Row 44: Guillaume / RC / 1.5h
$ws.Cells.Item(44, 2).Value = 42284
$ws.Cells.Item(44, 3).Value = "Guillaume"
$ws.Cells.Item(44, 4).Value = "RC"
$ws.Cells.Item(44, 6).Value = 1.5

# Row 45: Élodie / RC / 1.5h
$ws.Cells.Item(45, 2).Value = 42284
$ws.Cells.Item(45, 3).Value = "Élodie"
$ws.Cells.Item(45, 4).Value = "RC"
$ws.Cells.Item(45, 6).Value = 1.5

# Row 46: Olivier / RC / 1.5h
$ws.Cells.Item(46, 2).Value = 42284
$ws.Cells.Item(46, 3).Value = "Olivier"
$ws.Cells.Item(46, 4).Value = "RC"
$ws.Cells.Item(46, 6).Value = 1.5

# Row 47: Guillaume / AN / 1.5h
$ws.Cells.Item(47, 2).Value = 42284
$ws.Cells.Item(47, 3).Value = "Guillaume"
$ws.Cells.Item(47, 4).Value = "AN"
$ws.Cells.Item(47, 6).Value = 1.5

# Row 48: Élodie / BD / 1.5h
$ws.Cells.Item(48, 2).Value = 42284
$ws.Cells.Item(48, 3).Value = "Élodie"
$ws.Cells.Item(48, 4).Value = "BD"
$ws.Cells.Item(48, 6).Value = 1.5

# Row 49: Olivier / AN / 1.5h
$ws.Cells.Item(49, 2).Value = 42284
$ws.Cells.Item(49, 3).Value = "Olivier"
$ws.Cells.Item(49, 4).Value = "AN"
$ws.Cells.Item(49, 6).Value = 1.5

# --- Task descriptions (column E), written in original authoring order so
# brand-new shared-string entries land in the same order as the source edit.

# Row 25: existing entry, description extended with "+ DN"
$ws.Cells.Item(25, 5).Value = "DCO + DS + DN"

# Row 41: Progammation Employe
$ws.Cells.Item(41, 5).Value = "Progammation Employe"

# Row 42: Programmation Equipe (already used elsewhere on the sheet)
$ws.Cells.Item(42, 5).Value = "Programmation Equipe"

# Row 43: Documentation DCC Equipe
$ws.Cells.Item(43, 5).Value = "Documentation DCC Equipe"

# Row 44: 1ere Recontre pacakage test
$ws.Cells.Item(44, 5).Value = "1ere Recontre pacakage test"

# Row 45 & 46: 1ème Recontre pacakage test
$ws.Cells.Item(45, 5).Value = "1ème Recontre pacakage test"
$ws.Cells.Item(46, 5).Value = "1ème Recontre pacakage test"

# Row 48: Continuer la base de données Employe
$ws.Cells.Item(48, 5).Value = "Continuer la base de données Employe"

# Row 47 & 49: Prototype interface papier test
$ws.Cells.Item(47, 5).Value = "Prototype interface papier test"
$ws.Cells.Item(49, 5).Value = "Prototype interface papier test"

# --- Sheet view: scroll/select to follow the newly entered rows ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E53").Select()
